# Update EC database: swap the two worker records in rows 16 and 17
# (row 16 becomes ANDERSON MARRIAGA RODRIGUEZ, row 17 becomes MANUEL VICENTE ESPITIA SUAREZ)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1001970544"
$ws.Range("D16").Value = "ANDERSON MARRIAGA RODRIGUEZ"
$ws.Range("E16").Value = "1704"
$ws.Range("F16").Value = 25575

$ws.Range("C17").Value = "11059945"
$ws.Range("D17").Value = "MANUEL VICENTE ESPITIA SUAREZ"
$ws.Range("E17").Value = "1706"
$ws.Range("F17").Value = 5902
